$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '66.877.82'
$ws.Cells.Item(2, 5).Value = '  +2.09%  '
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '3.100.33'
$ws.Cells.Item(3, 5).Value = '  +5.24%  '
$ws.Cells.Item(4, 5).Value = '  -0.02%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '580.59'
$ws.Cells.Item(5, 5).Value = '  +1.70%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '172.93'
$ws.Cells.Item(6, 5).Value = '  +6.35%  '
$ws.Cells.Item(7, 5).Value = '  +0.06%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '3.094.90'
$ws.Cells.Item(8, 5).Value = '  +5.12%  '
$ws.Cells.Item(9, 5).Value = '  +1.49%  '
$ws.Cells.Item(10, 5).Value = '  -3.25%  '
$ws.Cells.Item(11, 5).Value = '  +3.81%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.483'
$ws.Cells.Item(12, 5).Value = '  +4.23%  '
$ws.Cells.Item(13, 5).Value = '  +2.48%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '37.39'
$ws.Cells.Item(14, 5).Value = '  +7.45%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '3.612.13'
$ws.Cells.Item(16, 5).Value = '  +5.06%  '
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '66.872.62'
$ws.Cells.Item(17, 5).Value = '  +1.99%  '
$ws.Cells.Item(18, 5).Value = '  +1.95%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '3.099.81'
$ws.Cells.Item(19, 5).Value = '  +5.18%  '
$ws.Cells.Item(20, 5).Value = '  +3.62%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '481.33'
$ws.Cells.Item(21, 5).Value = '  +8.17%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '0.716'
$ws.Cells.Item(22, 5).Value = '  +3.20%  '
$ws.Cells.Item(23, 5).Value = '  +3.41%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '83.99'
$ws.Cells.Item(24, 5).Value = '  +2.27%  '
$ws.Cells.Item(25, 5).Value = '  +4.72%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '13.15'
$ws.Cells.Item(26, 5).Value = '  +7.56%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '10.07'
$ws.Cells.Item(27, 5).Value = '  +0.61%  '
$ws.Cells.Item(28, 5).Value = '  +0.07%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '7.98'
$ws.Cells.Item(29, 5).Value = '  -0.70%  '
$ws.Cells.Item(30, 5).Value = '  -2.98%  '
$ws.Cells.Item(31, 5).Value = '  +3.70%  '
$ws.Cells.Item(32, 2).Value = 'EthereumClassic'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '28.79'
$ws.Cells.Item(32, 5).Value = '  +6.09%  '
$ws.Cells.Item(33, 2).Value = 'PEPE'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '0.0000100'
$ws.Cells.Item(33, 5).Value = '  -0.32%  '
$ws.Cells.Item(34, 5).Value = '  +0.99%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.999'
$ws.Cells.Item(35, 5).Value = '  -0.06%  '
$ws.Cells.Item(36, 5).Value = '  +3.35%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.993'
$ws.Cells.Item(37, 5).Value = '  +2.16%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '48.16'
$ws.Cells.Item(38, 5).Value = '  +4.03%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '2.12'
$ws.Cells.Item(39, 5).Value = '  +7.63%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '50.20'
$ws.Cells.Item(40, 5).Value = '  +2.32%  '
$ws.Cells.Item(41, 5).Value = '  +3.84%  '
$ws.Cells.Item(42, 5).Value = '  +0.67%  '
$ws.Cells.Item(43, 5).Value = '  +1.95%  '
$ws.Cells.Item(44, 5).Value = '  -0.34%  '
$ws.Cells.Item(45, 2).Value = 'Maker'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '2.840.96'
$ws.Cells.Item(45, 5).Value = '  +6.13%  '
$ws.Cells.Item(46, 2).Value = 'VeChain'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.0362'
$ws.Cells.Item(46, 5).Value = '  +3.02%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '385.67'
$ws.Cells.Item(47, 5).Value = '  +0.41%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '134.63'
$ws.Cells.Item(48, 5).Value = '  +0.78%  '
$ws.Cells.Item(49, 5).Value = '  -0.01%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '24.97'
$ws.Cells.Item(50, 5).Value = '  +4.24%  '
$ws.Cells.Item(51, 5).Value = '  +2.90%  '

$wb.Save()